# Updating barcode_offset and umi_offset
#
# 1. barcode_offset sheet: "0" -> "0,38,76" and append a new option "10,48,86"
# 2. umi_offset sheet: prepend a new option "1" (before "Not applicable", "16")
# 3. RNAseq sheet: widen the data-validation list ranges for columns O (barcode_offset)
#    and R (umi_offset) to include the newly added rows
# 4. .metadata sheet: bump pav:createdOn timestamp

$wb = $excel.ActiveWorkbook

# --- barcode_offset -------------------------------------------------------
$boSheet = $wb.Worksheets.Item("barcode_offset")
$boSheet.Range("A1").Value = "0,38,76"
$boSheet.Range("A5").NumberFormat = "@"
$boSheet.Range("A5").Value = "10,48,86"

# --- umi_offset -------------------------------------------------------------
$uoSheet = $wb.Worksheets.Item("umi_offset")
$uoSheet.Rows.Item(1).Insert()
$uoSheet.Range("A1").NumberFormat = "@"
$uoSheet.Range("A1").Value = "1"

# --- RNAseq data validation ranges ------------------------------------------
$main = $wb.Worksheets.Item("RNAseq")

$barcodeOffsetDv = $main.Range("O2:O1001")
$barcodeOffsetDv.Validation.Modify(3, 1, 1, "='barcode_offset'!`$A`$1:`$A`$5")
$barcodeOffsetDv.Validation.IgnoreBlank = $true
$barcodeOffsetDv.Validation.ShowError = $true
$barcodeOffsetDv.Validation.ErrorTitle = "Validation Error"
$barcodeOffsetDv.Validation.ErrorMessage = ""

$umiOffsetDv = $main.Range("R2:R1001")
$umiOffsetDv.Validation.Modify(3, 1, 1, "='umi_offset'!`$A`$1:`$A`$3")
$umiOffsetDv.Validation.IgnoreBlank = $true
$umiOffsetDv.Validation.ShowError = $true
$umiOffsetDv.Validation.ErrorTitle = "Validation Error"
$umiOffsetDv.Validation.ErrorMessage = ""

# --- .metadata: bump pav:createdOn ------------------------------------------
$metaSheet = $wb.Worksheets.Item(".metadata")
$metaSheet.Range("C2").Value = "2023-10-31T13:53:33-07:00"
